# Update the "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect the latest generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 3-6 and row 10
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 97
$wsExhibit.Range("F4").Value = 492
$wsExhibit.Range("F5").Value = 4878
$wsExhibit.Range("F6").Value = 375
$wsExhibit.Range("F10").Value = 229

# Sheet "全部类型": rows 3-6 and row 11
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 97
$wsAll.Range("F4").Value = 492
$wsAll.Range("F5").Value = 4878
$wsAll.Range("F6").Value = 375
$wsAll.Range("F11").Value = 229
